# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on the Kraken_Profits
# leve tables for several leves across the ALC/ARM/CRP/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9230.77
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 9230.77
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 9230.77
$ws.Range("N40").Value = -9580.77
$ws.Range("M40").ClearContents()

$ws.Range("H58").Value = 5000
$ws.Range("J58").Value = 5000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15300

$ws.Range("H80").Value = 287.16666
$ws.Range("J80").Value = 229.66667
$ws.Range("L80").Value = 689.00001
$ws.Range("N80").Value = -2685.00001

$ws.Range("H83").Value = 287.16666
$ws.Range("J83").Value = 229.66667
$ws.Range("L83").Value = 2067.00003
$ws.Range("N83").Value = -12051.00003

$ws.Range("H97").Value = 6775
$ws.Range("J97").Value = 6775
$ws.Range("L97").Value = 20325
$ws.Range("N97").Value = -21317

$ws.Range("H141").Value = 5825.1113
$ws.Range("I141").Value = 5825.1113
$ws.Range("K141").Value = 17475.3339
$ws.Range("M141").Value = -12295.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 673.5833
$ws.Range("I30").Value = 673.5833
$ws.Range("K30").Value = 673.5833
$ws.Range("M30").Value = -523.5833

$ws.Range("H32").Value = 3748
$ws.Range("I32").Value = 3748
$ws.Range("K32").Value = 3748
$ws.Range("M32").Value = -3461

$ws.Range("H35").Value = 2047
$ws.Range("I35").Value = 2047
$ws.Range("K35").Value = 2047
$ws.Range("M35").Value = -1641

$ws.Range("H36").Value = 2110.4
$ws.Range("I36").Value = 2110.4
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2110.4
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1764.4
$ws.Range("N36").ClearContents()

$ws.Range("H53").Value = 3900
$ws.Range("I53").Value = 3900
$ws.Range("K53").Value = 3900
$ws.Range("M53").Value = -3218

$ws.Range("H61").Value = 3080.3635
$ws.Range("I61").Value = 3188.4
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 3188.4
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -2976.4
$ws.Range("N61").Value = -2424

$ws.Range("H132").Value = 3209.25
$ws.Range("I132").Value = 3239.1428
$ws.Range("K132").Value = 9717.428400000001
$ws.Range("M132").Value = -7187.428400000001

$ws.Range("H136").Value = 3080.3635
$ws.Range("I136").Value = 3188.4
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 9565.200000000001
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -7015.200000000001
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 27433.334
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 27433.334
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 27433.334
$ws.Range("N50").Value = -28683.334
$ws.Range("M50").ClearContents()

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H59").Value = 20000
$ws.Range("I59").Value = 20000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 20000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -18855
$ws.Range("N59").ClearContents()

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 107.25
$ws.Range("I2").Value = 112.57143
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 112.57143
$ws.Range("L2").Value = 70
$ws.Range("M2").Value = 0.4285699999999935
$ws.Range("N2").Value = -296

$ws.Range("H15").Value = 30000
$ws.Range("J15").Value = 30000
$ws.Range("L15").Value = 30000
$ws.Range("N15").Value = -30576

$ws.Range("H80").Value = 6587.7
$ws.Range("I80").Value = 6234.75
$ws.Range("J80").Value = 7999.5
$ws.Range("K80").Value = 6234.75
$ws.Range("L80").Value = 7999.5
$ws.Range("M80").Value = -5236.75
$ws.Range("N80").Value = -9995.5

$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996

$ws.Range("H83").Value = 6587.7
$ws.Range("I83").Value = 6234.75
$ws.Range("J83").Value = 7999.5
$ws.Range("K83").Value = 31173.75
$ws.Range("L83").Value = 39997.5
$ws.Range("M83").Value = -26181.75
$ws.Range("N83").Value = -49981.5

$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984

$ws.Range("H97").Value = 912
$ws.Range("I97").Value = 916.3333
$ws.Range("K97").Value = 916.3333
$ws.Range("M97").Value = -420.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3390.1428
$ws.Range("I7").Value = 3649.1
$ws.Range("K7").Value = 3649.1
$ws.Range("M7").Value = -3537.1

$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 1000
$ws.Range("M46").Value = -812

$ws.Range("H55").Value = 6000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H126").Value = 3390.1428
$ws.Range("I126").Value = 3649.1
$ws.Range("K126").Value = 10947.3
$ws.Range("M126").Value = -8477.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 90000
$ws.Range("I75").Value = 90000
$ws.Range("K75").Value = 90000
$ws.Range("M75").Value = -89064

$ws.Range("H78").Value = 90000
$ws.Range("I78").Value = 90000
$ws.Range("K78").Value = 270000
$ws.Range("M78").Value = -265320

$ws.Range("H126").Value = 4999.75
$ws.Range("I126").Value = 4666.3335
$ws.Range("K126").Value = 13999.0005
$ws.Range("M126").Value = -11529.0005
